$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("labels")

# Fix mixed-up labels for Q7 (rows 46-56) and Q24 (rows 110-123)
# by restoring the correct Column C (full label) / Column D (short label) values.
$ws.Range("C46").Value = 'Access to highspeed internet'
$ws.Range("D46").Value = 'Internet access'
$ws.Range("C47").Value = 'Aging infrastructure (e.g., roads, levees, bridges)'
$ws.Range("D47").Value = 'Infrastructure'
$ws.Range("C48").Value = 'Climate change/Global warming'
$ws.Range("D48").Value = 'Climate change'
$ws.Range("C49").Value = 'Delta Conveyance/Delta Tunnel projects'
$ws.Range("D49").Value = 'Delta Conveyance'
$ws.Range("C50").Value = 'Environmental decline'
$ws.Range("D50").Value = 'Env decline'
$ws.Range("C51").Value = 'Lack of job or education opportunities'
$ws.Range("D51").Value = 'Job opportunities'
$ws.Range("C52").Value = 'Threats to public safety (e.g., first responders, trespassing, vandalis'
$ws.Range("D52").Value = 'Public safety'
$ws.Range("C53").Value = 'Social inequality'
$ws.Range("D53").Value = 'Social inequality'
$ws.Range("C54").Value = 'Traffic congestion/Limited transportation options'
$ws.Range("D54").Value = 'Transportation'
$ws.Range("C55").Value = 'Urban/Suburban development'
$ws.Range("D55").Value = 'Urb development'
$ws.Range("C56").Value = 'Other.'
$ws.Range("D56").Value = 'Other'
$ws.Range("C110").Value = 'Air filters'
$ws.Range("D110").Value = 'Air filters'
$ws.Range("C111").Value = 'Backup power supply/Generator'
$ws.Range("D111").Value = 'Backup power'
$ws.Range("C112").Value = 'Personal computer with internet connection'
$ws.Range("D112").Value = 'Computer internet'
$ws.Range("C113").Value = 'Mobile device with internet connection'
$ws.Range("D113").Value = 'Mobile internet'
$ws.Range("C114").Value = 'Domestic well for drinking water'
$ws.Range("D114").Value = 'Drinking well'
$ws.Range("C115").Value = 'Sewage system'
$ws.Range("D115").Value = 'Sewage'
$ws.Range("C116").Value = 'Homeowner''s or renter''s insurance'
$ws.Range("D116").Value = 'Home insurance'
$ws.Range("C117").Value = 'Flood insurance'
$ws.Range("D117").Value = 'Flood insurance'
$ws.Range("C118").Value = 'Earthquake insurance'
$ws.Range("D118").Value = 'Earthquake insurance'
$ws.Range("C119").Value = 'Health insurance'
$ws.Range("D119").Value = 'Health insurance'
$ws.Range("C120").Value = 'Personal motorized vehicle such as car, truck, motorbike, etc.'
$ws.Range("D120").Value = 'Personal vehicle'
$ws.Range("C121").Value = 'Public transportation such as bus or train route'
$ws.Range("D121").Value = 'Public transit'
$ws.Range("C122").Value = 'Emergency financial resources (e.g., savings, credit, loans)'
$ws.Range("D122").Value = 'Emergency finances'
$ws.Range("C123").Value = 'Family, friends, or supportive community you could stay with in the case of an emergency event evacuation'
$ws.Range("D123").Value = 'Emergency accomodations'

# Update the saved view/selection state of the sheet to match the author's edit.
$ws.Activate()
$win = $excel.ActiveWindow
try { $win.ScrollRow = 35 } catch {}
try { $win.ScrollColumn = 3 } catch {}
try { $win.TopLeftCell = $ws.Range("C35") } catch {}
$ws.Range("D63").Select()

